$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the match data (columns B:AC) between rows 127 and 128
$row127 = $ws.Range("B127:AC127").Value2
$row128 = $ws.Range("B128:AC128").Value2

$ws.Range("B127:AC127").Value2 = $row128
$ws.Range("B128:AC128").Value2 = $row127

# Swap the match data (columns B:AC) between rows 130 and 131
$row130 = $ws.Range("B130:AC130").Value2
$row131 = $ws.Range("B131:AC131").Value2

$ws.Range("B130:AC130").Value2 = $row131
$ws.Range("B131:AC131").Value2 = $row130
